# Apply WRI input data update to the Cash Flow Quantization Size workbook.

$wb = $excel.ActiveWorkbook

$wsCFQS = $wb.Worksheets.Item("CFQS")

# Update the label to clarify units (USD).
$wsCFQS.Range("B1").Value = "Quantization Size ($)"

# Update the quantization size value per WRI's first-draft 2.1 input data.
$wsCFQS.Range("B2").Value = 400000

# Widen column B to fit the new, longer label text (~19.43 chars).
$wsCFQS.Columns.Item(2).ColumnWidth = 18.67
